$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# B.Mayfield (row 2)
$rushing.Cells.Item(2, 5).Value = 8   # E2 (3DATT)
$rushing.Cells.Item(2, 6).Value = 1   # F2 (RZATT)

# N.Chubb (row 4)
$rushing.Cells.Item(4, 3).Value = 117 # C4 (1DATT)
$rushing.Cells.Item(4, 4).Value = 68  # D4 (2DATT)
$rushing.Cells.Item(4, 6).Value = 33  # F4 (RZATT)

# D.Johnson (row 6)
$rushing.Cells.Item(6, 3).Value = 11  # C6 (1DATT)
$rushing.Cells.Item(6, 5).Value = 5   # E6 (3DATT)

# J.Landry (row 11)
$rushing.Cells.Item(11, 4).Value = 3  # D11 (2DATT)

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# N.Chubb (row 2)
$receiving.Cells.Item(2, 3).Value = 16  # C2 (Short Target)

# D.Johnson (row 4)
$receiving.Cells.Item(4, 3).Value = 10  # C4 (Short Target)
$receiving.Cells.Item(4, 4).Value = 7   # D4 (Short Comp)

# J.Landry (row 8)
$receiving.Cells.Item(8, 3).Value = 72  # C8 (Short Target)
$receiving.Cells.Item(8, 4).Value = 53  # D8 (Short Comp)
$receiving.Cells.Item(8, 5).Value = 25  # E8 (Deep Target)
$receiving.Cells.Item(8, 6).Value = 13  # F8 (Deep Comp)
$receiving.Cells.Item(8, 7).Value = 9   # G8 (RZ Target)

# D.Peoples-Jones (row 9)
$receiving.Cells.Item(9, 3).Value = 25  # C9 (Short Target)
$receiving.Cells.Item(9, 4).Value = 17  # D9 (Short Comp)
$receiving.Cells.Item(9, 5).Value = 21  # E9 (Deep Target)
$receiving.Cells.Item(9, 6).Value = 10  # F9 (Deep Comp)

# R.Higgins (row 10)
$receiving.Cells.Item(10, 3).Value = 33 # C10 (Short Target)
$receiving.Cells.Item(10, 4).Value = 20 # D10 (Short Comp)
$receiving.Cells.Item(10, 5).Value = 9  # E10 (Deep Target)

# A.Schwartz (row 11)
$receiving.Cells.Item(11, 3).Value = 12 # C11 (Short Target)

# A.Hooper (row 13)
$receiving.Cells.Item(13, 3).Value = 56 # C13 (Short Target)
$receiving.Cells.Item(13, 4).Value = 34 # D13 (Short Comp)
$receiving.Cells.Item(13, 5).Value = 7  # E13 (Deep Target)
$receiving.Cells.Item(13, 6).Value = 4  # F13 (Deep Comp)

# D.Njoku (row 14)
$receiving.Cells.Item(14, 3).Value = 33 # C14 (Short Target)
$receiving.Cells.Item(14, 4).Value = 23 # D14 (Short Comp)
$receiving.Cells.Item(14, 7).Value = 7  # G14 (RZ Target)
$receiving.Cells.Item(14, 8).Value = 4  # H14 (RZ Comp)

# H.Bryant (row 15)
$receiving.Cells.Item(15, 3).Value = 17 # C15 (Short Target)
$receiving.Cells.Item(15, 4).Value = 14 # D15 (Short Comp)
$receiving.Cells.Item(15, 7).Value = 3  # G15 (RZ Target)
$receiving.Cells.Item(15, 8).Value = 3  # H15 (RZ Comp)
